$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Alfonso%Fasano%NULL%1,           Emanuele%Cereda%e.cereda@smatteo.pv.it%2,           Emanuele%Cereda%e.cereda@smatteo.pv.it%0,           Michela%Barichella%NULL%2,           Michela%Barichella%NULL%0,           Erica%Cassani%NULL%1,           Valentina%Ferri%NULL%1,           Anna Lena%Zecchinelli%NULL%1,           Gianni%Pezzoli%NULL%1]"
$ws.Range("E3").Value = "[Antonio%D’Avolio%NULL%1,           Valeria%Avataneo%NULL%2,           Valeria%Avataneo%NULL%0,           Alessandra%Manca%NULL%1,           Jessica%Cusato%NULL%1,           Amedeo%De Nicolò%NULL%2,           Amedeo%De Nicolò%NULL%0,           Renzo%Lucchini%NULL%1,           Franco%Keller%NULL%1,           Marco%Cantù%NULL%1]"
$ws.Range("E4").Value = "[Petre Cristian%Ilie%dr.cristianilie@gmail.com%1,           Simina%Stefanescu%NULL%1,           Lee%Smith%NULL%2]"
$ws.Range("E5").Value = "[ Claire E.%Hastie%null%1,            Daniel F.%Mackay%null%1,            Frederick%Ho%null%1,            Carlos A.%Celis-Morales%null%1,            Srinivasa Vittal%Katikireddi%null%1,            Claire L.%Niedzwiedz%null%1,            Bhautesh D.%Jani%null%1,            Paul%Welsh%null%1,            Frances S.%Mair%null%1,            Stuart R.%Gray%null%1,            Catherine A.%O’Donnell%null%1,            Jason MR.%Gill%null%1,            Naveed%Sattar%null%1,            Jill P.%Pell%null%1,          Claire E.%Hastie%null%1,          Daniel F.%Mackay%null%1,          Frederick%Ho%null%1,          Carlos A.%Celis-Morales%null%1,          Srinivasa Vittal%Katikireddi%null%1,          Claire L.%Niedzwiedz%null%1,          Bhautesh D.%Jani%null%1,          Paul%Welsh%null%0,          Frances S.%Mair%null%0,          Stuart R.%Gray%null%1,          Catherine A.%O’Donnell%null%0,          Jason MR.%Gill%null%1,          Naveed%Sattar%null%0,          Jill P.%Pell%null%0]"

Write-Host "Done updating E2:E5 author strings"
